# Replace "Create a basic" with "Edit" in the ALI: task bullet that reads
# "Create a basic index.html for our application further instructions will
# be given by skype." -> "Edit index.html for our application further
# instructions will be given by skype."

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Create a basic index.html",   # FindText
    $true,                          # MatchCase
    $false,                         # MatchWholeWord
    $false,                         # MatchWildcards
    $false,                         # MatchSoundsLike
    $false,                         # MatchAllWordForms
    $true,                          # Forward
    1,                               # Wrap (wdFindContinue)
    $false,                         # Format
    "Edit index.html",              # ReplaceWith
    2                                # Replace (wdReplaceAll)
)
